$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 137, pushing the existing
# rows 137:163 down to 139:165 (matches the diff: dimension grows from
# A1:T163 to A1:T165, two brand-new records appear at 137/138, and all
# previously-existing rows 137..163 retain their data but are renumbered
# 139..165).
$ws.Rows.Item(137).Insert()
$ws.Rows.Item(137).Insert()

# New row 137 - Primera, 2023-10-10
$ws.Cells.Item(137, 1).Value = 5
$ws.Cells.Item(137, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(137, 3).Value = "Maule"
$ws.Cells.Item(137, 4).Value = 45209
$ws.Cells.Item(137, 5).Value = 7
$ws.Cells.Item(137, 6).Value = "Fruta"
$ws.Cells.Item(137, 7).Value = 100107
$ws.Cells.Item(137, 8).Value = "Otros"
$ws.Cells.Item(137, 9).Value = 100107002
$ws.Cells.Item(137, 10).Value = "Chirimoya"
$ws.Cells.Item(137, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(137, 12).Value = "Primera"
$ws.Cells.Item(137, 13).Value = 340
$ws.Cells.Item(137, 14).Value = 19000
$ws.Cells.Item(137, 15).Value = 19000
$ws.Cells.Item(137, 16).Value = 19000
$ws.Cells.Item(137, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(137, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(137, 19).Value = 1900
$ws.Cells.Item(137, 20).Value = 10

# New row 138 - Segunda, 2023-10-10
$ws.Cells.Item(138, 1).Value = 5
$ws.Cells.Item(138, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(138, 3).Value = "Maule"
$ws.Cells.Item(138, 4).Value = 45209
$ws.Cells.Item(138, 5).Value = 7
$ws.Cells.Item(138, 6).Value = "Fruta"
$ws.Cells.Item(138, 7).Value = 100107
$ws.Cells.Item(138, 8).Value = "Otros"
$ws.Cells.Item(138, 9).Value = 100107002
$ws.Cells.Item(138, 10).Value = "Chirimoya"
$ws.Cells.Item(138, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(138, 12).Value = "Segunda"
$ws.Cells.Item(138, 13).Value = 280
$ws.Cells.Item(138, 14).Value = 15000
$ws.Cells.Item(138, 15).Value = 15000
$ws.Cells.Item(138, 16).Value = 15000
$ws.Cells.Item(138, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(138, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(138, 19).Value = 1500
$ws.Cells.Item(138, 20).Value = 10

# Make sure the date cells keep the expected date/time number format (same
# style used by the rest of column D).
$ws.Range("D137:D138").NumberFormat = "YYYY-MM-DD HH:MM:SS"
